$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data. D and E columns are forced to text
# (leading apostrophe) so that numeric-looking strings (e.g. "560.06", "0.995")
# are preserved exactly as text instead of being parsed into floating point numbers.

# Row 2
$ws.Range('D2').Value = '''61.858.21'
$ws.Range('E2').Value = '''  -0.80%  '

# Row 3
$ws.Range('D3').Value = '''2.398.96'

# Row 4
$ws.Range('E4').Value = '''  +0.00%  '

# Row 5
$ws.Range('D5').Value = '''560.06'
$ws.Range('E5').Value = '''  +0.58%  '

# Row 6
$ws.Range('D6').Value = '''142.22'
$ws.Range('E6').Value = '''  -1.12%  '

# Row 7
$ws.Range('E7').Value = '''  +0.09%  '

# Row 8
$ws.Range('D8').Value = '''0.532'
$ws.Range('E8').Value = '''  -0.86%  '

# Row 9
$ws.Range('E9').Value = '''  -1.71%  '

# Row 10
$ws.Range('E10').Value = '''  -1.88%  '

# Row 11
$ws.Range('E11').Value = '''  -2.97%  '

# Row 12
$ws.Range('E12').Value = '''  -1.12%  '

# Row 13
$ws.Range('D13').Value = '''25.47'
$ws.Range('E13').Value = '''  -3.13%  '

# Row 14
$ws.Range('E14').Value = '''  -1.92%  '

# Row 15
$ws.Range('D15').Value = '''2.831.81'
$ws.Range('E15').Value = '''  -1.20%  '

# Row 16
$ws.Range('D16').Value = '''61.748.52'
$ws.Range('E16').Value = '''  -0.85%  '

# Row 17
$ws.Range('D17').Value = '''2.403.38'
$ws.Range('E17').Value = '''  -0.96%  '

# Row 18
$ws.Range('E18').Value = '''  +0.62%  '

# Row 19
$ws.Range('D19').Value = '''320.20'
$ws.Range('E19').Value = '''  -1.52%  '

# Row 20
$ws.Range('E20').Value = '''  -1.13%  '

# Row 21
$ws.Range('E21').Value = '''  +0.32%  '

# Row 22
$ws.Range('E22').Value = '''  -0.17%  '

# Row 23
$ws.Range('D23').Value = '''65.51'
$ws.Range('E23').Value = '''  +0.75%  '

# Row 24
$ws.Range('E24').Value = '''  -2.94%  '

# Row 25
$ws.Range('E25').Value = '''  -4.72%  '

# Row 26
$ws.Range('D26').Value = '''560.67'
$ws.Range('E26').Value = '''  -2.26%  '

# Row 27
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '''2.517.66'
$ws.Range('E27').Value = '''  -0.95%  '

# Row 28
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '''0.995'
$ws.Range('E28').Value = '''  -0.23%  '

# Row 29
$ws.Range('E29').Value = '''  -2.77%  '

# Row 30
$ws.Range('D30').Value = '''8.17'
$ws.Range('E30').Value = '''  -2.97%  '

# Row 31
$ws.Range('E31').Value = '''  -5.17%  '

# Row 32
$ws.Range('E32').Value = '''  -1.58%  '

# Row 33
$ws.Range('E33').Value = '''  -0.21%  '

# Row 34
$ws.Range('E34').Value = '''  -4.81%  '

# Row 35
$ws.Range('E35').Value = '''  +0.09%  '

# Row 36
$ws.Range('D36').Value = '''4.72'
$ws.Range('E36').Value = '''  -2.27%  '

# Row 37
$ws.Range('D37').Value = '''152.40'
$ws.Range('E37').Value = '''  +2.72%  '

# Row 38
$ws.Range('E38').Value = '''  -6.01%  '

# Row 39
$ws.Range('E39').Value = '''  -2.33%  '

# Row 40
$ws.Range('E40').Value = '''  -1.84%  '

# Row 41
$ws.Range('E41').Value = '''  -6.09%  '

# Row 42
$ws.Range('E42').Value = '''  -0.04%  '

# Row 43
$ws.Range('E43').Value = '''  -3.73%  '

# Row 44
$ws.Range('D44').Value = '''147.15'
$ws.Range('E44').Value = '''  -3.39%  '

# Row 45
$ws.Range('D45').Value = '''3.60'
$ws.Range('E45').Value = '''  -1.24%  '

# Row 46
$ws.Range('E46').Value = '''  -3.13%  '

# Row 47
$ws.Range('D47').Value = '''19.75'
$ws.Range('E47').Value = '''  -3.60%  '

# Row 48
$ws.Range('E48').Value = '''  -0.97%  '

# Row 49
$ws.Range('E49').Value = '''  +0.16%  '

# Row 50
$ws.Range('E50').Value = '''  -2.02%  '
